$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 57.357143
$ws.Range("I6").Value = 49.5
$ws.Range("J6").Value = 104.5
$ws.Range("K6").Value = 148.5
$ws.Range("L6").Value = 313.5
$ws.Range("M6").Value = -36.5
$ws.Range("N6").Value = -537.5

$ws.Range("H12").Value = 891
$ws.Range("I12").Value = 891
$ws.Range("K12").Value = 891
$ws.Range("M12").Value = -721

$ws.Range("H40").Value = 3508.818
$ws.Range("I40").Value = 4257.143
$ws.Range("K40").Value = 4257.143
$ws.Range("M40").Value = -4082.143

$ws.Range("H46").Value = 0
$ws.Range("I46").Value = 0
$ws.Range("J46").Value = 0
$ws.Range("K46").Value = 0
$ws.Range("L46").Value = 0
$ws.Range("M46").ClearContents()
$ws.Range("N46").ClearContents()

$ws.Range("H52").Value = 501.66666
$ws.Range("I52").Value = 501.66666
$ws.Range("K52").Value = 1504.99998
$ws.Range("M52").Value = -1344.99998

$ws.Range("H60").Value = 0
$ws.Range("I60").Value = 0
$ws.Range("J60").Value = 0
$ws.Range("K60").Value = 0
$ws.Range("L60").Value = 0
$ws.Range("M60").ClearContents()
$ws.Range("N60").ClearContents()

$ws.Range("H98").Value = 2000
$ws.Range("I98").Value = 2000
$ws.Range("K98").Value = 2000
$ws.Range("M98").Value = -502

$ws.Range("H116").Value = 4967
$ws.Range("I116").Value = 4967
$ws.Range("K116").Value = 4967
$ws.Range("M116").Value = -1525

$ws.Range("H122").Value = 2000
$ws.Range("I122").Value = 2000
$ws.Range("K122").Value = 6000
$ws.Range("M122").Value = -3550

$ws.Range("H141").Value = 446
$ws.Range("I141").Value = 446
$ws.Range("K141").Value = 1338
$ws.Range("M141").Value = 3842

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H30").Value = 1325
$ws.Range("I30").Value = 1325
$ws.Range("K30").Value = 1325
$ws.Range("M30").Value = -1175

$ws.Range("H32").Value = 4855.8125
$ws.Range("I32").Value = 4855.8125
$ws.Range("K32").Value = 4855.8125
$ws.Range("M32").Value = -4568.8125

$ws.Range("H34").Value = 18997
$ws.Range("J34").Value = 18997
$ws.Range("L34").Value = 18997
$ws.Range("N34").Value = -19539

$ws.Range("H37").Value = 0
$ws.Range("I37").Value = 0
$ws.Range("J37").Value = 0
$ws.Range("K37").Value = 0
$ws.Range("L37").Value = 0
$ws.Range("M37").ClearContents()
$ws.Range("N37").ClearContents()

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H45").Value = 35000
$ws.Range("J45").Value = 35000
$ws.Range("L45").Value = 35000
$ws.Range("N45").Value = -36616

$ws.Range("H105").Value = 1207.5714
$ws.Range("I105").Value = 1207.5714
$ws.Range("K105").Value = 1207.5714
$ws.Range("M105").Value = 539.4286

$ws.Range("H134").Value = 2969.2856
$ws.Range("I134").Value = 2660
$ws.Range("J134").Value = 3742.5
$ws.Range("K134").Value = 7980
$ws.Range("L134").Value = 11227.5
$ws.Range("M134").Value = -5445
$ws.Range("N134").Value = -16297.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H10").Value = 1883
$ws.Range("I10").Value = 1259.6
$ws.Range("J10").Value = 5000
$ws.Range("K10").Value = 1259.6
$ws.Range("L10").Value = 5000
$ws.Range("M10").Value = -1120.6
$ws.Range("N10").Value = -5278

$ws.Range("H16").Value = 864.3333
$ws.Range("I16").Value = 899
$ws.Range("K16").Value = 899
$ws.Range("M16").Value = -612

$ws.Range("H56").Value = 46000
$ws.Range("I56").Value = 0
$ws.Range("J56").Value = 46000
$ws.Range("K56").Value = 0
$ws.Range("L56").Value = 46000
$ws.Range("M56").ClearContents()
$ws.Range("N56").Value = -47690

$ws.Range("H99").Value = 26499.75
$ws.Range("I99").Value = 34666.668
$ws.Range("K99").Value = 34666.668
$ws.Range("M99").Value = -33168.668

$ws.Range("H113").Value = 864.3333
$ws.Range("I113").Value = 899
$ws.Range("K113").Value = 899
$ws.Range("M113").Value = 1271

$ws.Range("H122").Value = 6227
$ws.Range("I122").Value = 878.6667
$ws.Range("K122").Value = 2636.0001
$ws.Range("M122").Value = -186.0001000000002

$ws.Range("H126").Value = 26499.75
$ws.Range("I126").Value = 34666.668
$ws.Range("K126").Value = 104000.004
$ws.Range("M126").Value = -101530.004

$ws.Range("H134").Value = 2996
$ws.Range("J134").Value = 3000
$ws.Range("L134").Value = 9000
$ws.Range("N134").Value = -14070

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H38").Value = 156.57143
$ws.Range("I38").Value = 113
$ws.Range("J38").Value = 174
$ws.Range("K38").Value = 339
$ws.Range("L38").Value = 522
$ws.Range("M38").Value = 8
$ws.Range("N38").Value = -1216

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H7").Value = 168.5
$ws.Range("I7").Value = 2
$ws.Range("J7").Value = 2000
$ws.Range("K7").Value = 2
$ws.Range("L7").Value = 2000
$ws.Range("M7").Value = 110
$ws.Range("N7").Value = -2224

$ws.Range("H8").Value = 168.5
$ws.Range("I8").Value = 2
$ws.Range("J8").Value = 2000
$ws.Range("K8").Value = 2
$ws.Range("L8").Value = 2000
$ws.Range("M8").Value = 137
$ws.Range("N8").Value = -2278

$ws.Range("H20").Value = 0
$ws.Range("I20").Value = 0
$ws.Range("K20").Value = 0
$ws.Range("M20").ClearContents()

$ws.Range("H22").Value = 4953.375
$ws.Range("I22").Value = 1127
$ws.Range("J22").Value = 5500
$ws.Range("K22").Value = 1127
$ws.Range("L22").Value = 5500
$ws.Range("M22").Value = -598
$ws.Range("N22").Value = -6558

$ws.Range("H24").Value = 0
$ws.Range("J24").Value = 0
$ws.Range("L24").Value = 0
$ws.Range("N24").ClearContents()

$ws.Range("H97").Value = 338
$ws.Range("I97").Value = 338
$ws.Range("K97").Value = 338
$ws.Range("M97").Value = 158

$ws.Range("H102").Value = 1000.1429
$ws.Range("I102").Value = 1040.2
$ws.Range("K102").Value = 1040.2
$ws.Range("M102").Value = 581.8

$ws.Range("H113").Value = 992.875
$ws.Range("I113").Value = 992.875
$ws.Range("J113").Value = 0
$ws.Range("K113").Value = 992.875
$ws.Range("L113").Value = 0
$ws.Range("M113").Value = 1177.125
$ws.Range("N113").ClearContents()

$ws.Range("H122").Value = 4198.8887
$ws.Range("I122").Value = 2760.3333
$ws.Range("K122").Value = 8280.999899999999
$ws.Range("M122").Value = -5830.999899999999

$ws.Range("H126").Value = 2999
$ws.Range("I126").Value = 2999
$ws.Range("K126").Value = 8997
$ws.Range("M126").Value = -6527

$ws.Range("H132").Value = 2512
$ws.Range("I132").Value = 2011.3334
$ws.Range("K132").Value = 6034.0002
$ws.Range("M132").Value = -3504.0002

$ws.Range("H136").Value = 0
$ws.Range("J136").Value = 0
$ws.Range("L136").Value = 0
$ws.Range("N136").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 1322.625
$ws.Range("I16").Value = 1028
$ws.Range("J16").Value = 1499.4
$ws.Range("K16").Value = 1028
$ws.Range("L16").Value = 1499.4
$ws.Range("M16").Value = -858
$ws.Range("N16").Value = -1839.4

$ws.Range("H31").Value = 15026.25
$ws.Range("I31").Value = 105
$ws.Range("K31").Value = 105
$ws.Range("M31").Value = 143

$ws.Range("H46").Value = 1500
$ws.Range("I46").Value = 1500
$ws.Range("K46").Value = 1500
$ws.Range("M46").Value = -1312

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H4").Value = 10002
$ws.Range("I4").Value = 10002
$ws.Range("J4").Value = 0
$ws.Range("K4").Value = 10002
$ws.Range("L4").Value = 0
$ws.Range("M4").Value = -9889
$ws.Range("N4").ClearContents()

$ws.Range("H28").Value = 19000
$ws.Range("J28").Value = 19000
$ws.Range("L28").Value = 19000
$ws.Range("N28").Value = -19696

$ws.Range("H100").Value = 178.5
$ws.Range("J100").Value = 0
$ws.Range("L100").Value = 0
$ws.Range("N100").ClearContents()

$ws.Range("H107").Value = 518.8125
$ws.Range("I107").Value = 391.75
$ws.Range("K107").Value = 1175.25
$ws.Range("M107").Value = 744.75

$ws.Range("H113").Value = 336.85715
$ws.Range("I113").Value = 267.25
$ws.Range("J113").Value = 429.66666
$ws.Range("K113").Value = 801.75
$ws.Range("L113").Value = 1288.99998
$ws.Range("M113").Value = 1368.25
$ws.Range("N113").Value = -5628.999980000001

$ws.Range("H123").Value = 40429
$ws.Range("J123").Value = 40429
$ws.Range("L123").Value = 40429
$ws.Range("N123").Value = -50229

$ws.Range("H125").Value = 0
$ws.Range("J125").Value = 0
$ws.Range("L125").Value = 0
$ws.Range("N125").ClearContents()

$ws.Range("H126").Value = 3491.5334
$ws.Range("I126").Value = 2631.0833
$ws.Range("J126").Value = 6933.3335
$ws.Range("K126").Value = 7893.249899999999
$ws.Range("L126").Value = 20800.0005
$ws.Range("M126").Value = -5423.249899999999
$ws.Range("N126").Value = -25740.0005

$ws.Range("H132").Value = 9999
$ws.Range("I132").Value = 9999
$ws.Range("K132").Value = 29997
$ws.Range("M132").Value = -27467
